$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" values (column D) look like plain numbers to Excel
# (e.g. "589.83"), but in this sheet the Price column stores text values.
# Force those specific cells to Text format before writing so Excel keeps them
# as strings instead of auto-converting them to numeric values, then restore the
# default (Normal) cell style once the text has been written.
$textCells = @("D5", "D6", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D41", "D42", "D45", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.622.40"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "3.188.80"
$ws.Range("E3").Value = "  -3.28%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "589.83"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("D6").Value = "136.26"
$ws.Range("E6").Value = "  -4.03%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.183.85"
$ws.Range("E8").Value = "  -3.40%  "
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  -4.57%  "
$ws.Range("D11").Value = "5.28"
$ws.Range("E11").Value = "  -3.65%  "
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").Value = "0.0000236"
$ws.Range("E13").Value = "  -4.43%  "
$ws.Range("D14").Value = "33.41"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "3.714.14"
$ws.Range("E15").Value = "  -3.27%  "
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "3.191.61"
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("D18").Value = "62.639.10"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "6.54"
$ws.Range("E19").Value = "  -4.74%  "
$ws.Range("D20").Value = "457.72"
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").Value = "13.94"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "0.704"
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("D23").Value = "7.64"
$ws.Range("E23").Value = "  -4.65%  "
$ws.Range("D24").Value = "13.35"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "83.61"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "6.88"
$ws.Range("E29").Value = "  -6.20%  "
$ws.Range("D30").Value = "7.79"
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("D31").Value = "2.02"
$ws.Range("E31").Value = "  -6.42%  "
$ws.Range("D32").Value = "27.41"
$ws.Range("E32").Value = "  -5.77%  "
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").Value = "2.39"
$ws.Range("E34").Value = "  -5.41%  "
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "51.14"
$ws.Range("E37").Value = "  -3.28%  "
$ws.Range("E38").Value = "  -5.88%  "
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "402.51"
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D42").Value = "8.02"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").Value = "2.841.81"
$ws.Range("E43").Value = "  -6.88%  "
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").Value = "36.60"
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("E46").Value = "  -5.91%  "
$ws.Range("D47").Value = "2.14"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "125.24"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").Value = "25.60"
$ws.Range("E50").Value = "  -2.48%  "
$ws.Range("E51").Value = "  -3.44%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
